$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property -> set literal text value "false"
# (use a formula that evaluates to the text, then paste-as-values so the
# result lands as plain shared-string text instead of being auto-coerced
# to a Boolean by a direct .Value assignment)
$ws.Range("B7").Formula = "=""false"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

# Row 14 = "Case Sensitive" property -> set literal text value "true"
$ws.Range("B14").Formula = "=""true"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)

$excel.CutCopyMode = 0
